$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the opening sentence was split across three runs (the
# original docx-generation tool broke it at the "#18" placeholder).
# Re-join it into a single run with the full sentence. We overwrite the
# spanning range's text twice: once with a scratch placeholder (so the
# engine sees a genuine change and actually rewrites the range instead
# of treating an identical re-assignment as a no-op), then with the
# real final text. This collapses the three backing runs into one
# (their formatting is identical) while keeping straight apostrophes
# (setting Range.Text does not go through AutoCorrect's smart-quote
# substitution the way Find.Execute's replacement text does).
# ---------------------------------------------------------------------
$sentence = "This Deed of Lease is made at #1 this #2 day of #3, #18 between #4 of #5 hereinafter called 'The Lessor' of the One Part and #6 also of #7 hereinafter called 'The Lessee' of the Other Part."

$full = $d.Content.Text
$start = $full.IndexOf("This Deed of Lease is made at #1 this #2 day of #3, ")
$oldLen = "This Deed of Lease is made at #1 this #2 day of #3, #18 between #4 of #5 hereinafter called 'The Lessor' of the One Part and #6 also of #7 hereinafter called 'The Lessee' of the Other Part.".Length
$rng = $d.Range($start, $start + $oldLen)
$rng.Text = "TEMP_PLACEHOLDER_8E1F2A"
$full2 = $d.Content.Text
$start2 = $full2.IndexOf("TEMP_PLACEHOLDER_8E1F2A")
$rng2 = $d.Range($start2, $start2 + "TEMP_PLACEHOLDER_8E1F2A".Length)
$rng2.Text = $sentence

# ---------------------------------------------------------------------
# Change 2: "Withinnamed Lessee #7 in the presence of #17" should read
# #4 instead of #7 (the placeholder was wrong), and the diff shows the
# result keeps the placeholder in its own run (split out from the
# surrounding text), matching how the original doc's other placeholders
# are laid out. Wrapping the target text in a bookmark before editing
# forces the engine to keep the edited span as its own run instead of
# re-coalescing it with identically-formatted neighbours; removing the
# bookmark afterwards leaves no trace of it behind.
# ---------------------------------------------------------------------
$full3 = $d.Content.Text
$idx7 = $full3.IndexOf("#7 in the presence of #17")
$target = $d.Range($idx7, $idx7 + 2)
$bm = $d.Bookmarks.Add("tmp_fix_bm", $target)
$target.Text = "#4"
$d.Bookmarks("tmp_fix_bm").Delete()

Write-Output "Applied both fixes."
